$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 177 in column C (Fitness) need to be updated to 7293
$ws.Range("C2:C177").Value = 7293
